$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.109.66"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "1.801.07"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'316.87"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.5443"
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("D8").Value = "'0.3791"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.07481"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").Value = "'41.97"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "'1.096"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'6.237"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'20.47"
$ws.Range("D15").Value = "'7.376"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "1.792.66"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "'89.45"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'0.06542"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'17.42"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'5.944"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "28.135.28"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").Value = "'156.01"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").Value = "'20.41"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "2.006.97"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "'2.328"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "'121.89"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").Value = "'0.1109"
$ws.Range("E31").Value = "  +7.92%  "
$ws.Range("D32").Value = "'1.116"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").Value = "'3.670"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").Value = "'5.556"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'0.06880"
$ws.Range("E35").Value = "  +6.19%  "
$ws.Range("D36").Value = "'0.2218"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("D37").Value = "'0.02296"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").Value = "'5.086"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").Value = "'8.449"
$ws.Range("E39").Value = "  -5.68%  "
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").Value = "'0.6165"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("D42").Value = "'1.175"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "'1.420"
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("D44").Value = "'13.30"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "'3.685"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "'0.5749"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").Value = "'1.920"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").Value = "'0.06818"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").Value = "'0.00000000299"
$ws.Range("E51").Value = "  +39.94%  "
